$d = $word.ActiveDocument

# 1. "the different categories of sandwiches" -> "discovering different type of sandwiches"
$d.Content.Find.Execute("he is very interested in the different categories of sandwiches", $true, $false, $false, $false, $false, $true, 1, $false, "he is very interested in discovering different type of sandwiches", 2)

# 2. "startup" -> "start-up"
$d.Content.Find.Execute("Andrew has created his own startup", $true, $false, $false, $false, $false, $true, 1, $false, "Andrew has created his own start-up", 2)

# 3. "the nutritional value of his meals" -> "the calories included in his meals"
$d.Content.Find.Execute("His primary concern is the nutritional value of his meals", $true, $false, $false, $false, $false, $true, 1, $false, "His primary concern is the calories included in his meals", 2)

# 4. Move the _GoBack bookmark to right after "discovering different type"
$r = $d.Content
$r.Find.Execute("discovering different type", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmRange = $r.Duplicate
$bmRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRange)
